$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 47: ChangeRelationCardinality change operator (first parameter line)
# Set the brand-new shared strings in the same order the target workbook
# introduces them (CHANGECARDINALITY, then ChangeRelationCardinality, then
# the cardinality parameter description) so sharedStrings.xml ends up with
# the same ordering as the authored edit.
$ws.Cells.Item(47, 4).Value = "CHANGECARDINALITY"          # D47
$ws.Cells.Item(47, 2).Value = "ChangeRelationCardinality"  # B47
$ws.Cells.Item(47, 3).Value = "Relation"                   # C47
$ws.Cells.Item(47, 5).Value = "relation"                   # E47
$ws.Cells.Item(47, 6).Value = "Relation"                   # F47

# New row 48: second parameter line (the new cardinality value)
$ws.Cells.Item(48, 5).Value = "cardinality (newCardinality?)"  # E48
$ws.Cells.Item(48, 6).Value = "Cardinality"                     # F48

# Reflect the new selection/scroll position left behind by the edit.
$ws.Range("E49").Select()

Write-Host "Added ChangeRelationCardinality rows 47-48"
